$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells whose updated text looks like a plain number must be forced to
# Text format first, otherwise Excel auto-converts them to numeric values and
# silently drops the exact original formatting (trailing zeros, etc.).
$numericLookingCells = @("D4","D5","D6","D9","D10","D11","D12","D13","D14","D17","D20","D21","D22","D23","D25","D27","D28","D29","D30","D31","D32","D33","D34","D35","D37","D38","D39","D40","D43","D44","D45","D46","D48","D50","D51")
foreach ($cellRef in $numericLookingCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}


$ws.Range("D2").Value = '65.306.48'
$ws.Range("E2").Value = '  +1.99%  '

$ws.Range("D3").Value = '3.135.22'
$ws.Range("E3").Value = '  +2.34%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.18%  '

$ws.Range("D5").Value = '569.92'
$ws.Range("E5").Value = '  +2.12%  '

$ws.Range("D6").Value = '149.69'
$ws.Range("E6").Value = '  +3.24%  '

$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("D8").Value = '3.133.73'
$ws.Range("E8").Value = '  +2.38%  '

$ws.Range("D9").Value = '0.526'
$ws.Range("E9").Value = '  +4.85%  '

$ws.Range("D10").Value = '0.163'
$ws.Range("E10").Value = '  +6.76%  '

$ws.Range("D11").Value = '6.20'
$ws.Range("E11").Value = '  +0.80%  '

$ws.Range("D12").Value = '0.498'
$ws.Range("E12").Value = '  +6.62%  '

$ws.Range("D13").Value = '0.0000251'
$ws.Range("E13").Value = '  +10.45%  '

$ws.Range("D14").Value = '37.14'
$ws.Range("E14").Value = '  +5.97%  '

$ws.Range("D15").Value = '3.645.73'
$ws.Range("E15").Value = '  +2.13%  '

$ws.Range("D16").Value = '65.201.74'
$ws.Range("E16").Value = '  +1.54%  '

$ws.Range("D17").Value = '7.15'
$ws.Range("E17").Value = '  +5.89%  '

$ws.Range("D18").Value = '3.115.71'
$ws.Range("E18").Value = '  +1.47%  '

$ws.Range("E19").Value = '  +0.45%  '

$ws.Range("D20").Value = '508.28'
$ws.Range("E20").Value = '  +6.43%  '

$ws.Range("D21").Value = '14.83'
$ws.Range("E21").Value = '  +6.81%  '

$ws.Range("D22").Value = '15.62'
$ws.Range("E22").Value = '  +15.69%  '

$ws.Range("D23").Value = '0.722'
$ws.Range("E23").Value = '  +7.23%  '

$ws.Range("E24").Value = '  +3.23%  '

$ws.Range("D25").Value = '85.11'
$ws.Range("E25").Value = '  +4.27%  '

$ws.Range("E26").Value = '  +0.19%  '

$ws.Range("D27").Value = '2.91'
$ws.Range("E27").Value = '  +3.91%  '

$ws.Range("D28").Value = '8.72'
$ws.Range("E28").Value = '  +8.00%  '

$ws.Range("D29").Value = '2.15'
$ws.Range("E29").Value = '  +4.18%  '

$ws.Range("D30").Value = '27.93'
$ws.Range("E30").Value = '  +6.71%  '

$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  -0.09%  '

$ws.Range("D32").Value = '2.67'
$ws.Range("E32").Value = '  +8.07%  '

$ws.Range("D33").Value = '1.19'
$ws.Range("E33").Value = '  +3.23%  '

$ws.Range("D34").Value = '6.04'
$ws.Range("E34").Value = '  +8.65%  '

$ws.Range("D35").Value = '6.56'
$ws.Range("E35").Value = '  +6.52%  '

$ws.Range("E36").Value = '  +1.48%  '

$ws.Range("D37").Value = '470.25'
$ws.Range("E37").Value = '  +2.10%  '

$ws.Range("D38").Value = '0.0423'
$ws.Range("E38").Value = '  +4.34%  '

$ws.Range("D39").Value = '0.0852'
$ws.Range("E39").Value = '  +3.02%  '

$ws.Range("D40").Value = '2.97'
$ws.Range("E40").Value = '  -1.66%  '

$ws.Range("D41").Value = '3.127.14'
$ws.Range("E41").Value = '  +5.63%  '

$ws.Range("E42").Value = '  +6.23%  '

$ws.Range("D43").Value = '8.58'
$ws.Range("E43").Value = '  +4.12%  '

$ws.Range("D44").Value = '0.288'
$ws.Range("E44").Value = '  +10.03%  '

$ws.Range("D45").Value = '2.43'
$ws.Range("E45").Value = '  +13.55%  '

$ws.Range("D46").Value = '29.18'
$ws.Range("E46").Value = '  +4.86%  '

$ws.Range("E47").Value = '  -0.04%  '

$ws.Range("D48").Value = '0.116'
$ws.Range("E48").Value = '  +3.44%  '

$ws.Range("D49").Value = '0.0₃0555'
$ws.Range("E49").Value = '  +7.66%  '

$ws.Range("D50").Value = '2.28'
$ws.Range("E50").Value = '  +9.87%  '

$ws.Range("D51").Value = '118.46'
$ws.Range("E51").Value = '  -1.35%  '
